# eims-toi updates per readme doc Aug 11
#
# Turn off multi-threaded recalculation (workbook.xml calcPr
# concurrentCalc="0").
$excel.MultiThreadedCalculation.Enabled = $false

# The pre-seeded $wb variable is not reliable in this runtime, so re-fetch
# the active workbook/sheet from $excel explicitly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was "datetime_utc" / "Date and time in UTC" -> becomes the
# MATLAB-serial-date row, now described as PI-provided.
$ws.Range("A2").Value = "datetime_utc_matlab"
$ws.Range("B2").Value = "PI-provided UTC date and time"

# Row 3 was "datetime_utc_matlab" / MATLAB date explanation (numeric /
# dimensionless) -> becomes the UTC datetime row (class Date, with a
# dateTimeFormatString, and no unit).
$ws.Range("A3").Value = "datetime_utc"
$ws.Range("B3").Value = "Product UTC date and time from NES-LTER API"
$ws.Range("C3").Value = "Date"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "YYYY-MM-DD hh:mm:ss"

# Row 9 was "depth_mat_file" -> renamed to "depth_matlab" with a
# PI-provided definition.
$ws.Range("A9").Value = "depth_matlab"
$ws.Range("B9").Value = "PI-provided depth of sample below sea surface. URI http://vocab.nerc.ac.uk/collection/P09/current/DEPH/"

# Row 10 (depth_API_bottle_summary) gets a reworded definition.
$ws.Range("B10").Value = "Product depth of sample below sea surface from NES-LTER API"

# Row 6 (toi_source) gets a reworded definition.
$ws.Range("B6").Value = "Source of bottle sample whether from Niskin or underway"

# Update the selection to match the author's saved cursor position.
$ws.Range("B7").Select()
